# Índice de Frequência Ampliado
# Adds a new "calcular_indices_ampliados" function block to the
# Funcoes_Inputs and Funcoes_Outputs sheets, and updates the cell
# selections that were left behind on a few sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Funcoes_Inputs: append the input rows for calcular_indices_ampliados
# ---------------------------------------------------------------------
$wsInputs = $wb.Worksheets.Item("Funcoes_Inputs")

$inputRows = @(
    @("Funcionarios", $true),
    @("Nev_Safast_Tipico", $false),
    @("Nev_Safast_Trajeto", $false),
    @("Nev_Safast_DoenOcup", $false),
    @("Nev_Safast_NRelac", $false),
    @("Nev_Obito_Tipico", $false),
    @("Nev_Obito_Trajeto", $false),
    @("Nev_Obito_DoenOcup", $false),
    @("Nev_Obito_NRelac", $false),
    @("Nev_Afmenor15_Tipico", $false),
    @("Nev_Afmenor15_Trajeto", $false),
    @("Nev_Afmenor15_DoenOcup", $false),
    @("Nev_Afmenor15_NRelac", $false),
    @("Nev_Afmaior15_Tipico", $false),
    @("Nev_Afmaior15_Trajeto", $false),
    @("Nev_Afmaior15_DoenOcup", $false),
    @("Nev_Afmaior15_NRelac", $false)
)

$startRow = 97
for ($i = 0; $i -lt $inputRows.Count; $i++) {
    $r = $startRow + $i
    $name = $inputRows[$i][0]
    $bold = $inputRows[$i][1]

    $wsInputs.Cells.Item($r, 1).Value = "calcular_indices_ampliados"
    $wsInputs.Cells.Item($r, 2).Value = $name
    if ($bold) {
        $wsInputs.Cells.Item($r, 2).Font.Bold = $true
    }
}

# ---------------------------------------------------------------------
# 2. Funcoes_Outputs: append the output rows for calcular_indices_ampliados
# ---------------------------------------------------------------------
$wsOutputs = $wb.Worksheets.Item("Funcoes_Outputs")

$wsOutputs.Cells.Item(49, 1).Value = "calcular_indices_ampliados"
$wsOutputs.Cells.Item(49, 2).Value = "EventosIndiceFrequenciaAmpliado"

$wsOutputs.Cells.Item(50, 1).Value = "calcular_indices_ampliados"
$wsOutputs.Cells.Item(50, 2).Value = "IndiceFrequenciaAmpliado"

# ---------------------------------------------------------------------
# 3. Restore / update the leftover cell selections on each sheet
# ---------------------------------------------------------------------
$wsDados = $wb.Worksheets.Item("Dados_Projetados")
[void]$wsDados.Range("E14").Select()

$wsParametros = $wb.Worksheets.Item("Parametros")
[void]$wsParametros.Range("A1").Select()

[void]$wsInputs.Range("B98").Select()

# Funcoes_Outputs must stay the active sheet/tab, so select its range last
# and make sure it ends up as the active sheet.
[void]$wsOutputs.Activate()
[void]$wsOutputs.Range("B47").Select()
